# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff/Handback DateTime"
# timestamps that get refreshed each time the handback report is regenerated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$overview.Range("G2").Value = "2016-08-29 17:13:48"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2) for the first file row.
$zhcn.Range("H2").Value = "2016-08-29 17:13:42"
$zhcn.Range("K2").Value = "2016-08-29 17:13:59"

# de-de sheet: "Correspond Handback DateTime" (K2) for the first file row.
$dede.Range("K2").Value = "2016-08-29 17:14:15"
